$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Wed Feb 15 12:49:21 EST 2023"
$ws.Range("B3").Value = "Wed Feb 15 12:50:04 EST 2023"
$ws.Range("B4").Value = "Wed Feb 15 12:50:51 EST 2023"
$ws.Range("B5").Value = "Wed Feb 15 12:51:33 EST 2023"
$ws.Range("B6").Value = "Wed Feb 15 12:52:14 EST 2023"
$ws.Range("B7").Value = "Wed Feb 15 12:52:55 EST 2023"
$ws.Range("B8").Value = "Wed Feb 15 12:53:35 EST 2023"
$ws.Range("B9").Value = "Wed Feb 15 12:54:22 EST 2023"
$ws.Range("B10").Value = "Wed Feb 15 12:55:12 EST 2023"
$ws.Range("B11").Value = "Wed Feb 15 12:55:59 EST 2023"
$ws.Range("B12").Value = "Wed Feb 15 12:56:42 EST 2023"
$ws.Range("B13").Value = "Wed Feb 15 13:01:45 EST 2023"
$ws.Range("B14").Value = "Wed Feb 15 13:02:29 EST 2023"
$ws.Range("B15").Value = "Wed Feb 15 13:03:14 EST 2023"
$ws.Range("B16").Value = "Wed Feb 15 13:04:00 EST 2023"
$ws.Range("B17").Value = "Wed Feb 15 13:04:45 EST 2023"
$ws.Range("B18").Value = "Wed Feb 15 13:05:32 EST 2023"
$ws.Range("B19").Value = "Wed Feb 15 13:06:22 EST 2023"
$ws.Range("B20").Value = "Wed Feb 15 13:07:17 EST 2023"
$ws.Range("B21").Value = "Wed Feb 15 13:08:11 EST 2023"
$ws.Range("B22").Value = "Wed Feb 15 13:09:02 EST 2023"

$ws.Range("A12").Value = "Fail"
